# Auto-generated edit script: the underlying observation records were
# re-sorted/re-exported upstream, so each sheet row (114-127) now shows a
# different observation's data while staying at the same row number.
# Columns A (Id) / B (Taxonsorteringsordning) / D (Rodlistade) / E
# (TaxonId) / F (Artnamn) / G (Vetenskapligt namn) / H (Auktor) / Q (Ost)
# / R (Nord) change together per row, and the optional columns M
# (Aktivitet) and AJ/AK/AO (Substrat fields) are only set for the
# observations that actually carry that data - so some rows gain those
# cells (set) and others lose them (ClearContents).
#
# Old row -> new row each record's data now lives at:
#   114->124   117->126   120->122   123->116   126->115
#   115->123   118->118   121->121   124->127   127->125
#   116->114   119->117   122->120   125->119
# (values below are the literal target cell contents, taken directly
# from the target OOXML rather than re-derived/copied at runtime)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 114
$ws.Cells.Item(114, 1).Value = 111743517
$ws.Cells.Item(114, 2).Value = 73634
$ws.Cells.Item(114, 4).Value = "LC"
$ws.Cells.Item(114, 5).Value = 6426
$ws.Cells.Item(114, 6).Value = "Kattfotslav"
$ws.Cells.Item(114, 7).Value = "Felipes leucopellaeus"
$ws.Cells.Item(114, 8).Value = "(Ach.) Frisch & G.Thor"
$ws.Cells.Item(114, 17).Value = 339278.3213300391
$ws.Cells.Item(114, 18).Value = 6571107.378548244
$ws.Cells.Item(114, 36).ClearContents()
$ws.Cells.Item(114, 37).ClearContents()
$ws.Cells.Item(114, 41).ClearContents()

# Row 115
$ws.Cells.Item(115, 1).Value = 111743526
$ws.Cells.Item(115, 2).Value = 90666
$ws.Cells.Item(115, 4).Value = "LC"
$ws.Cells.Item(115, 5).Value = 4364
$ws.Cells.Item(115, 6).Value = "Dropptaggsvamp"
$ws.Cells.Item(115, 7).Value = "Hydnellum ferrugineum"
$ws.Cells.Item(115, 8).Value = "(Fr.:Fr.) P. Karst."
$ws.Cells.Item(115, 17).Value = 338870.1217119552
$ws.Cells.Item(115, 18).Value = 6571086.774471543

# Row 116
$ws.Cells.Item(116, 1).Value = 111743527
$ws.Cells.Item(116, 2).Value = 96348
$ws.Cells.Item(116, 4).Value = "VU"
$ws.Cells.Item(116, 5).Value = 220787
$ws.Cells.Item(116, 6).Value = "Knärot"
$ws.Cells.Item(116, 7).Value = "Goodyera repens"
$ws.Cells.Item(116, 8).Value = "(L.) R. Br."
$ws.Cells.Item(116, 17).Value = 338598.1684531783
$ws.Cells.Item(116, 18).Value = 6571109.585305012

# Row 117
$ws.Cells.Item(117, 1).Value = 111743519
$ws.Cells.Item(117, 2).Value = 90666
$ws.Cells.Item(117, 5).Value = 4364
$ws.Cells.Item(117, 6).Value = "Dropptaggsvamp"
$ws.Cells.Item(117, 7).Value = "Hydnellum ferrugineum"
$ws.Cells.Item(117, 8).Value = "(Fr.:Fr.) P. Karst."
$ws.Cells.Item(117, 17).Value = 339118.4126724883
$ws.Cells.Item(117, 18).Value = 6571062.424656671

# Row 119
$ws.Cells.Item(119, 1).Value = 111743520
$ws.Cells.Item(119, 2).Value = 56398
$ws.Cells.Item(119, 4).Value = "NT"
$ws.Cells.Item(119, 5).Value = 100109
$ws.Cells.Item(119, 6).Value = "Tretåig hackspett"
$ws.Cells.Item(119, 7).Value = "Picoides tridactylus"
$ws.Cells.Item(119, 8).Value = "(Linnaeus, 1758)"
$ws.Cells.Item(119, 13).Value = "färska spår"
$ws.Cells.Item(119, 17).Value = 339096.8530521042
$ws.Cells.Item(119, 18).Value = 6571013.66294401
$ws.Cells.Item(119, 36).Value = "gran"
$ws.Cells.Item(119, 37).Value = "Picea abies"
$ws.Cells.Item(119, 41).Value = "Picea abies"

# Row 120
$ws.Cells.Item(120, 1).Value = 111743546
$ws.Cells.Item(120, 17).Value = 339474.5644867857
$ws.Cells.Item(120, 18).Value = 6571113.931964876

# Row 122
$ws.Cells.Item(122, 1).Value = 111743551
$ws.Cells.Item(122, 17).Value = 339522.8608171764
$ws.Cells.Item(122, 18).Value = 6571091.407599592

# Row 123
$ws.Cells.Item(123, 1).Value = 111743554
$ws.Cells.Item(123, 2).Value = 88966
$ws.Cells.Item(123, 4).Value = "NT"
$ws.Cells.Item(123, 5).Value = 5754
$ws.Cells.Item(123, 6).Value = "Gultoppig fingersvamp"
$ws.Cells.Item(123, 7).Value = "Ramaria testaceoflava"
$ws.Cells.Item(123, 8).Value = "(Bres.) Corner"
$ws.Cells.Item(123, 17).Value = 339577.2032005055
$ws.Cells.Item(123, 18).Value = 6571127.007499221

# Row 124
$ws.Cells.Item(124, 1).Value = 111743524
$ws.Cells.Item(124, 2).Value = 94134
$ws.Cells.Item(124, 4).Value = "NT"
$ws.Cells.Item(124, 5).Value = 53
$ws.Cells.Item(124, 6).Value = "Vedtrappmossa"
$ws.Cells.Item(124, 7).Value = "Crossocalyx hellerianus"
$ws.Cells.Item(124, 8).Value = "(Nees ex Lindenb.) Meyl."
$ws.Cells.Item(124, 17).Value = 338949.7235384365
$ws.Cells.Item(124, 18).Value = 6571040.381812023
$ws.Cells.Item(124, 36).Value = "tall"
$ws.Cells.Item(124, 37).Value = "Pinus sylvestris"
$ws.Cells.Item(124, 41).Value = "Pinus sylvestris"

# Row 125
$ws.Cells.Item(125, 1).Value = 111743521
$ws.Cells.Item(125, 2).Value = 96348
$ws.Cells.Item(125, 4).Value = "VU"
$ws.Cells.Item(125, 5).Value = 220787
$ws.Cells.Item(125, 6).Value = "Knärot"
$ws.Cells.Item(125, 7).Value = "Goodyera repens"
$ws.Cells.Item(125, 8).Value = "(L.) R. Br."
$ws.Cells.Item(125, 17).Value = 339070.1946752003
$ws.Cells.Item(125, 18).Value = 6571001.989220584
$ws.Cells.Item(125, 13).ClearContents()
$ws.Cells.Item(125, 36).ClearContents()
$ws.Cells.Item(125, 37).ClearContents()
$ws.Cells.Item(125, 41).ClearContents()

# Row 126
$ws.Cells.Item(126, 1).Value = 111743523
$ws.Cells.Item(126, 2).Value = 73634
$ws.Cells.Item(126, 5).Value = 6426
$ws.Cells.Item(126, 6).Value = "Kattfotslav"
$ws.Cells.Item(126, 7).Value = "Felipes leucopellaeus"
$ws.Cells.Item(126, 8).Value = "(Ach.) Frisch & G.Thor"
$ws.Cells.Item(126, 17).Value = 339009.0243061834
$ws.Cells.Item(126, 18).Value = 6571011.238422027

# Row 127
$ws.Cells.Item(127, 1).Value = 111743515
$ws.Cells.Item(127, 17).Value = 339441.7613444271
$ws.Cells.Item(127, 18).Value = 6571017.506567059
